# "Colocando header nos gráficos"
# Adds a header label to column A (row 1) on each data sheet, fixes the
# accentuation of several Portuguese labels, removes the bold/border style
# from those relabeled cells (keeping it only on the new header row),
# drops the obsolete "Teto" row on the Emissoes sheet, and refreshes the
# cost figures (with a header) on the last sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share the same row layout (Hidro, Gas Natural, Carvao, ...)
# ---------------------------------------------------------------------
$rowLabels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

foreach ($i in 1..4) {
    $ws = $wb.Worksheets.Item($i)

    # New header cell for the technology/source column, styled like the
    # rest of row 1 (bold, bordered, centered) by copying B1's format.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    foreach ($r in $rowLabels.Keys) {
        $cell = $ws.Cells.Item($r, 1)
        $cell.Value = $rowLabels[$r]
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------
# Sheet 5 - "Emissoes Totais (MtCO2eq)": relabel, drop the "Teto" row
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws5.Cells.Item(2, 1).Value = "P.Médio"
$ws5.Cells.Item(2, 1).Style = "Normal"

$ws5.Cells.Item(3, 1).Value = "P.Crítico"
$ws5.Cells.Item(3, 1).Style = "Normal"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6 - "Custo Total (bilhões de R$)": new header + updated figures
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B1 becomes the text "2015" (same as the B1 header on the other sheets),
# keeping its existing bold/border style - paste-values from sheet 1's B1
# so the text stays a real string instead of Excel auto-typing it numeric.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws6.Cells.Item(2, 1).Value = "Expansão Centralizada"
$ws6.Cells.Item(2, 1).Style = "Normal"
$ws6.Cells.Item(2, 2).Value = 610

$ws6.Cells.Item(3, 1).Value = "Expansão por GD"
$ws6.Cells.Item(3, 1).Style = "Normal"
$ws6.Cells.Item(3, 2).Value = 67
